$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$countryCol = $ws.Range("A4:A216")

# --- Update country statistics with the latest COVID-19 numbers ---

# Estados Unidos
$r = $countryCol.Find("Estados Unidos").Row
$ws.Cells.Item($r, 2).Value = 654343
$ws.Cells.Item($r, 3).Value = 6195
$ws.Cells.Item($r, 5).Value = 564235
$ws.Cells.Item($r, 6).Value = 13369
$ws.Cells.Item($r, 7).Value = 902
$ws.Cells.Item($r, 8).Value = 33490

# Brasil
$r = $countryCol.Find("Brasil").Row
$ws.Cells.Item($r, 2).Value = 29214
$ws.Cells.Item($r, 3).Value = 604
$ws.Cells.Item($r, 5).Value = 13419
$ws.Cells.Item($r, 7).Value = 12
$ws.Cells.Item($r, 8).Value = 1769

# Austria
$r = $countryCol.Find("Austria").Row
$ws.Cells.Item($r, 2).Value = 14474
$ws.Cells.Item($r, 3).Value = 124
$ws.Cells.Item($r, 5).Value = 5095

# Mali's case count jumps above El Salvador's and Martinica's, so once the
# table (kept sorted by "Casos totales" descending) is re-sorted below, Mali
# moves two spots up the ranking.
$r = $countryCol.Find("Mali").Row
$ws.Cells.Item($r, 2).Value = 171
$ws.Cells.Item($r, 3).Value = 23
$ws.Cells.Item($r, 5).Value = 124

# --- Re-sort the country data range (A4:H216) by column B ("Casos totales"), descending ---
$dataRange = $ws.Range("A4:H216")
$sortKey = $ws.Range("B4:B216")
$dataRange.Sort($sortKey, 2, $null, $null, 1, $null, 1, 2, $false, 2, $null, 1)
